$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the "Setting up the location/area of service ... Could be
#    unnecessary and add extra complexity)" bullet entirely (including its
#    paragraph mark), which merges the following bullet
#    ("Accepting/Refusing delivery orders") up into its place.
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "Setting up the location/area of service") {
        $startPara = $p
    }
    if ($startPara -ne $null -and $t -match "Accepting/Refusing delivery orders") {
        $endPara = $p
        break
    }
}

$delStart = $startPara.Range.Start
$delEnd = $endPara.Range.Start
$d.Range($delStart, $delEnd).Delete()

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: it currently sits on the
#    "Order Available for Pickup" run; after the edit above it should sit on
#    the "Order Cancellation & Refund handling" run instead.
# ---------------------------------------------------------------------------
$pickupPara = $null
$cancelPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "Order Available for Pickup") {
        $pickupPara = $p
    }
    if ($t -match "Order Cancellation & Refund handling") {
        $cancelPara = $p
    }
}

# Drop the marker from "Order Available for Pickup" by reassigning its text
# (re-setting a range's Text clears any rendering-cache markers attached to
# it while leaving the visible content identical).
$pickupRange = $pickupPara.Range
$pickupRange.MoveEnd(1, -1) | Out-Null
$pickupRange.Text = $pickupRange.Text

# Add the marker onto "Order Cancellation & Refund handling" by round
# tripping its WordOpenXML with <w:lastRenderedPageBreak/> spliced in right
# before the run's <w:t>.
$cancelRange = $cancelPara.Range
$cancelRange.MoveEnd(1, -1) | Out-Null
$xml = $cancelRange.WordOpenXML()
$marker = "<w:lastRenderedPageBreak/>"
$needle = "<w:t xml:space=`"preserve`">Order Cancellation &amp; Refund handling "
$replacement = $marker + $needle
$xml2 = $xml.Replace($needle, $replacement)
$cancelRange.WordOpenXML = $xml2
